$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2024-04-03 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-04 Thursday", 2)

# Update the arithmetic table cells, row by row (1-based row indices with data: 1,5,9,13,17)
$tbl = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)
$values = @(
    @("60÷3=", "22÷9=", "79÷3=", "66÷2=", "75÷8="),
    @("34÷6=", "69÷6=", "39÷6=", "64÷3=", "33÷5="),
    @("61÷3=", "58÷3=", "76÷6=", "31÷8=", "28÷5="),
    @("75÷7=", "23÷7=", "95÷6=", "22÷3=", "30÷6="),
    @("71÷2=", "37÷3=", "69÷2=", "43÷5=", "54÷8=")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowIndex = $rows[$r]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $tbl.Cell($rowIndex, $c)
        $cell.Range.Text = $values[$r][$c - 1]
    }
}
